$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = 468
for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    $cell.Value2 = $cell.Value2 + 1
}
